# Apply updated crypto price / volume data (and a few reordered rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.934.28"
$ws.Range("E2").Value = "  +5.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.581.69"
$ws.Range("E3").Value = "  +5.92%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.59"
$ws.Range("E5").Value = "  +3.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.65"
$ws.Range("E6").Value = "  +11.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  +2.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.190"
$ws.Range("E9").Value = "  +12.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.580.06"
$ws.Range("E10").Value = "  +5.91%  "

$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  +6.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.68"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.112.41"
$ws.Range("E14").Value = "  +7.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.980.79"
$ws.Range("E15").Value = "  +6.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("E16").Value = "  +3.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.58"
$ws.Range("E17").Value = "  +10.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.591.14"
$ws.Range("E18").Value = "  +6.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.67"
$ws.Range("E19").Value = "  +25.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("E20").Value = "  +10.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.62"
$ws.Range("E21").Value = "  +8.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.18"
$ws.Range("E22").Value = "  +13.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.02"
$ws.Range("E23").Value = "  +4.87%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.85"
$ws.Range("E25").Value = "  +2.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.05"
$ws.Range("E26").Value = "  +9.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("E27").Value = "  +10.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.722.15"
$ws.Range("E28").Value = "  +6.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0908"
$ws.Range("E30").Value = "  +11.21%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.35"
$ws.Range("E31").Value = "  +16.09%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "491.14"
$ws.Range("E32").Value = "  +14.94%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.59"
$ws.Range("E33").Value = "  +6.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.71"
$ws.Range("E34").Value = "  +6.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +11.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.07"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.23"
$ws.Range("E38").Value = "  +1.16%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.87"
$ws.Range("E39").Value = "  +5.01%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.77"
$ws.Range("E41").Value = "  +10.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.63"
$ws.Range("E42").Value = "  +8.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.319"
$ws.Range("E43").Value = "  +6.92%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "154.01"
$ws.Range("E44").Value = "  +18.51%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0859"
$ws.Range("E45").Value = "  +19.43%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.49"
$ws.Range("E46").Value = "  +2.90%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.14"
$ws.Range("E47").Value = "  +6.33%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.27"
$ws.Range("E48").Value = "  +11.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.55"
$ws.Range("E49").Value = "  +6.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.516"
$ws.Range("E50").Value = "  +7.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.88"
$ws.Range("E51").Value = "  +17.99%  "

